$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Address Flow 1"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = "first flow 1"
$ws.Range("D3").Value = "last flow 2"
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = "City Flow 1"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "111"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1111"
